$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.593.09'
$ws.Range("E2").Value = '  -0.90%  '

$ws.Range("D3").Value = '1.864.38'
$ws.Range("E3").Value = '  -1.23%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.015'
$ws.Range("E4").Value = '  -0.25%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.17'
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.014'
$ws.Range("E6").Value = '  -0.28%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4651'
$ws.Range("E7").Value = '  -0.64%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3916'
$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.02'
$ws.Range("E9").Value = '  -3.36%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07967'
$ws.Range("E10").Value = '  -1.02%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.000'
$ws.Range("E11").Value = '  -1.52%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.67'
$ws.Range("E12").Value = '  -1.14%  '

$ws.Range("D13").Value = '1.863.79'
$ws.Range("E13").Value = '  -2.61%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.944'
$ws.Range("E14").Value = '  -0.29%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.207'
$ws.Range("E15").Value = '  +1.75%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.26'
$ws.Range("E17").Value = '  +1.11%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06730'
$ws.Range("E18").Value = '  -0.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001043'
$ws.Range("E19").Value = '  -0.53%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.19'
$ws.Range("E20").Value = '  -0.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.012'
$ws.Range("E21").Value = '  -0.40%  '

$ws.Range("D22").Value = '27.609.85'
$ws.Range("E22").Value = '  -0.90%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.455'
$ws.Range("E23").Value = '  -0.85%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.93'
$ws.Range("E24").Value = '  -0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.306'
$ws.Range("E25").Value = '  -1.58%  '

$ws.Range("D26").Value = '2.088.32'
$ws.Range("E26").Value = '  -2.38%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.43'
$ws.Range("E27").Value = '  -0.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.63'
$ws.Range("E28").Value = '  -2.25%  '

$ws.Range("E29").Value = '  +2.73%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.422'
$ws.Range("E30").Value = '  -0.74%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.62'
$ws.Range("E31").Value = '  -0.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9730'
$ws.Range("E32").Value = '  +0.74%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09442'
$ws.Range("E33").Value = '  -0.34%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.624'
$ws.Range("E34").Value = '  -0.64%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.308'
$ws.Range("E35").Value = '  -0.77%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.337'
$ws.Range("E36").Value = '  -5.26%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02232'
$ws.Range("E37").Value = '  -0.96%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06024'
$ws.Range("E38").Value = '  -1.63%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.313'
$ws.Range("E39").Value = '  +3.75%  '

$ws.Range("E40").Value = '  -1.54%  '

$ws.Range("E41").Value = '  -0.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5941'
$ws.Range("E42").Value = '  -0.73%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1867'
$ws.Range("E43").Value = '  -0.89%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.29'
$ws.Range("E44").Value = '  +0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.249'
$ws.Range("E45").Value = '  -1.30%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5583'
$ws.Range("E46").Value = '  -1.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.20'
$ws.Range("E47").Value = '  +0.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.922'
$ws.Range("E48").Value = '  -0.09%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06728'
$ws.Range("E49").Value = '  -2.85%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.71'
$ws.Range("E50").Value = '  -1.77%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.050'
$ws.Range("E51").Value = '  -1.86%  '
